# Applies the commit: inserts a new weekly price record row at row 125
# (pushing all subsequent rows down by one), extending the data range
# from A1:R175 to A1:R176.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before existing row 125, shifting rows 125-175 -> 126-176
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new record's data.
$ws.Cells.Item(125, 1).Value()  = 7
$ws.Cells.Item(125, 2).Value()  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(125, 3).Value()  = "Ñuble"
$ws.Cells.Item(125, 4).Value()  = "2023-07-12"
$ws.Cells.Item(125, 5).Value()  = 16
$ws.Cells.Item(125, 6).Value()  = 100112037
$ws.Cells.Item(125, 7).Value()  = "Cebollín"
$ws.Cells.Item(125, 8).Value()  = "Sin especificar"
$ws.Cells.Item(125, 9).Value()  = "Primera"
$ws.Cells.Item(125, 10).Value() = 200
$ws.Cells.Item(125, 11).Value() = 7000
$ws.Cells.Item(125, 12).Value() = 7000
$ws.Cells.Item(125, 13).Value() = 7000
$ws.Cells.Item(125, 14).Value() = "$/paquete 36 unidades"
$ws.Cells.Item(125, 15).Value() = "Provincia de Diguillín"
$ws.Cells.Item(125, 16).Value() = 194
$ws.Cells.Item(125, 17).Value() = 36
$ws.Cells.Item(125, 18).Value() = "Hortaliza"

$ws.Cells.Item(125, 4).NumberFormat() = "YYYY-MM-DD HH:MM:SS"
